$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple pairwise row-content swaps (columns B:AC) - preserves column A (sequence id)
$v1 = $ws.Range("B2:AC2").Value()
$v2 = $ws.Range("B3:AC3").Value()
$ws.Range("B2:AC2").Value = $v2
$ws.Range("B3:AC3").Value = $v1

$v1 = $ws.Range("B14:AC14").Value()
$v2 = $ws.Range("B15:AC15").Value()
$ws.Range("B14:AC14").Value = $v2
$ws.Range("B15:AC15").Value = $v1

$v1 = $ws.Range("B17:AC17").Value()
$v2 = $ws.Range("B18:AC18").Value()
$ws.Range("B17:AC17").Value = $v2
$ws.Range("B18:AC18").Value = $v1

$v1 = $ws.Range("B24:AC24").Value()
$v2 = $ws.Range("B25:AC25").Value()
$ws.Range("B24:AC24").Value = $v2
$ws.Range("B25:AC25").Value = $v1

$v1 = $ws.Range("B35:AC35").Value()
$v2 = $ws.Range("B36:AC36").Value()
$ws.Range("B35:AC35").Value = $v2
$ws.Range("B36:AC36").Value = $v1

$v1 = $ws.Range("B46:AC46").Value()
$v2 = $ws.Range("B47:AC47").Value()
$ws.Range("B46:AC46").Value = $v2
$ws.Range("B47:AC47").Value = $v1

$v1 = $ws.Range("B80:AC80").Value()
$v2 = $ws.Range("B81:AC81").Value()
$ws.Range("B80:AC80").Value = $v2
$ws.Range("B81:AC81").Value = $v1

$v1 = $ws.Range("B100:AC100").Value()
$v2 = $ws.Range("B101:AC101").Value()
$ws.Range("B100:AC100").Value = $v2
$ws.Range("B101:AC101").Value = $v1

$v1 = $ws.Range("B127:AC127").Value()
$v2 = $ws.Range("B128:AC128").Value()
$ws.Range("B127:AC127").Value = $v2
$ws.Range("B128:AC128").Value = $v1

$v1 = $ws.Range("B152:AC152").Value()
$v2 = $ws.Range("B153:AC153").Value()
$ws.Range("B152:AC152").Value = $v2
$ws.Range("B153:AC153").Value = $v1

$v1 = $ws.Range("B155:AC155").Value()
$v2 = $ws.Range("B156:AC156").Value()
$ws.Range("B155:AC155").Value = $v2
$ws.Range("B156:AC156").Value = $v1

# Rotation among rows 142-147 (content shifts up by one, wraps around)
$row142 = $ws.Range("B142:AC142").Value()
$row143 = $ws.Range("B143:AC143").Value()
$row144 = $ws.Range("B144:AC144").Value()
$row145 = $ws.Range("B145:AC145").Value()
$row146 = $ws.Range("B146:AC146").Value()
$row147 = $ws.Range("B147:AC147").Value()

$ws.Range("B142:AC142").Value = $row143
$ws.Range("B143:AC143").Value = $row144
$ws.Range("B144:AC144").Value = $row145
$ws.Range("B145:AC145").Value = $row146
$ws.Range("B146:AC146").Value = $row147
$ws.Range("B147:AC147").Value = $row142

# Remove the trailing placeholder fixtures (rows 157-159) with no result data yet
$ws.Rows("157:159").Delete()
